$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1), columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting of the existing header cells (bold, bordered, centered)
# by copying the format from H1 (an existing header cell) onto I1:J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data cells for columns I and J, rows 2-5
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 9
